$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# The newest entry (currently the last row of the "France" table) needs to move up
# so the table stays ordered the same way it was before the new row got appended at
# the bottom. Re-sort just the affected tail of the table (rows 105-116) by the
# first column so the freshly added player lands back in its proper slot while every
# other row in the table is left completely untouched.
$sortRange = $ws.Range("A105:D116")
$keyRange = $ws.Range("A105:A116")

$srt = $lo.Sort
$srt.SortFields.Clear()
$srt.SetRange($sortRange)
$srt.SortFields.Add($keyRange, $null, 1)
$srt.Header = 0
$srt.Apply()

Write-Host "Reordered rows 105-116 of table" $lo.Name
